$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.139.43'
$ws.Range('E2').Value = '  -3.30%  '
$ws.Range('D3').Value = '1.925.68'
$ws.Range('E3').Value = '  -2.42%  '
$ws.Range('E4').Value = '  -1.07%  '
$ws.Range('D5').Value = "'330.48"
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D7').Value = "'0.4724"
$ws.Range('E7').Value = '  -4.87%  '
$ws.Range('D8').Value = "'0.4055"
$ws.Range('E8').Value = '  -3.66%  '
$ws.Range('D9').Value = "'53.09"
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('D10').Value = "'0.08424"
$ws.Range('E10').Value = '  -9.17%  '
$ws.Range('E11').Value = '  -4.67%  '
$ws.Range('E12').Value = '  -2.51%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = "'7.511"
$ws.Range('E13').Value = '  -4.89%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'6.100"
$ws.Range('E14').Value = '  -5.51%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.853.99'
$ws.Range('E15').Value = '  -6.65%  '
$ws.Range('D16').Value = "'1.001"
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').Value = "'90.46"
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('E18').Value = '  -3.77%  '
$ws.Range('D19').Value = "'0.06582"
$ws.Range('E19').Value = '  -1.80%  '
$ws.Range('D20').Value = "'18.09"
$ws.Range('E20').Value = '  -5.56%  '
$ws.Range('D21').Value = "'1.001"
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').Value = "'5.748"
$ws.Range('E22').Value = '  -3.60%  '
$ws.Range('D23').Value = '28.130.28'
$ws.Range('E23').Value = '  -3.38%  '
$ws.Range('D24').Value = "'11.39"
$ws.Range('E24').Value = '  -4.78%  '
$ws.Range('D25').Value = "'2.286"
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').Value = '2.137.63'
$ws.Range('E26').Value = '  -3.76%  '
$ws.Range('D27').Value = "'154.28"
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('D28').Value = "'20.08"
$ws.Range('E28').Value = '  -3.16%  '
$ws.Range('D29').Value = "'2.149"
$ws.Range('E29').Value = '  -5.04%  '
$ws.Range('D30').Value = "'5.751"
$ws.Range('E30').Value = '  -8.28%  '
$ws.Range('D31').Value = "'123.74"
$ws.Range('E31').Value = '  -2.75%  '
$ws.Range('D32').Value = "'0.9777"
$ws.Range('E32').Value = '  -6.57%  '
$ws.Range('D33').Value = "'0.09615"
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('D34').Value = "'1.454"
$ws.Range('E34').Value = '  -3.65%  '
$ws.Range('D36').Value = "'3.637"
$ws.Range('E36').Value = '  -2.65%  '
$ws.Range('D37').Value = "'8.986"
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('E38').Value = '  -4.36%  '
$ws.Range('D39').Value = "'0.06176"
$ws.Range('E39').Value = '  -3.94%  '
$ws.Range('D40').Value = "'1.236"
$ws.Range('E40').Value = '  -6.89%  '
$ws.Range('E41').Value = '  -4.78%  '
$ws.Range('D42').Value = "'11.07"
$ws.Range('E42').Value = '  -3.90%  '
$ws.Range('E43').Value = '  -0.96%  '
$ws.Range('D44').Value = "'0.1905"
$ws.Range('E44').Value = '  -4.87%  '
$ws.Range('D45').Value = "'1.307"
$ws.Range('E45').Value = '  -4.37%  '
$ws.Range('D46').Value = "'0.5885"
$ws.Range('E46').Value = '  -5.18%  '
$ws.Range('D47').Value = "'12.78"
$ws.Range('E47').Value = '  -3.84%  '
$ws.Range('D48').Value = "'2.037"
$ws.Range('E48').Value = '  -6.73%  '
$ws.Range('D49').Value = "'3.473"
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').Value = "'0.06835"
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('D51').Value = "'110.30"
$ws.Range('E51').Value = '  -2.71%  '
